{"js": "// The cover page had a stray duplicate line reading \"of the\" (left over\n// from restructuring the \"... HND of / INFORMATION TECHNOLOGY / of the\"\n// block) \u2014 remove that extraneous paragraph entirely, exactly like a\n// person selecting the whole line (including its paragraph mark) and\n// pressing Delete.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"of the\") {\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# The cover page had a stray duplicate line reading \"of the\" (left over\n# from restructuring the \"... HND of / INFORMATION TECHNOLOGY / of the\"\n# block) \u2014 remove that extraneous paragraph entirely, exactly like a\n# person selecting the whole line (including its paragraph mark) and\n# pressing Delete.\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"of the\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
